# This workbook contains a weekly price log. A new weekly record was
# inserted as row 7 (pushing the former rows 7-43 down to rows 8-44).
# The new row reuses the same Calidad/Precio/Unidad/Origen data as the
# old row 8 record (which, after the shift below, now lives in row 9),
# but carries a new date (2022-04-26, serial 44677) and a new
# "Volumen" (120).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 7; rows 7:43 shift down to 8:44.
$ws.Rows("7:7").Insert()

# After the insert, row 9 holds the data that used to be in row 8
# (same Calidad/Precio/Unidad/Origen) - exactly the template needed
# for the brand-new row 7. Copy it over.
$ws.Range("A9:T9").Copy()
$ws.Range("A7:T7").PasteSpecial()

# Now overwrite the two fields that differentiate the new record.
$ws.Range("D7").Value = 44677
$ws.Range("M7").Value = 120
